$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row: "_old" -> "_FV2404", "_new" -> "_FV2410" ---
# Columns A..J carry the "_old" variants, column K is the literal "diff"
# column (unchanged), columns L..U carry the "_new" variants.
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($rightCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# --- 2) Turn the data range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3) Freeze the header row (row 1) ---
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Table:" $tbl.Name $tbl.Range.Address()
